# Apply the "add item art resource and item can random push to game ui" change:
#  - Add a new "UIName" column (F) to the WeaponData sheet, pushing the old
#    "Description" column from F to G.
#  - Split the old combined UIPath value
#    "Assets/ArtResources/Weapons/Weapons Sprite Sheet.png[Weapons Sprite Sheet_123]"
#    into a plain UIPath (E4) and the new UIName (F4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WeaponData")

# --- Row 2 / Row 3 / Row 4: move the existing "Description" column (F) data
#     into the new column G *before* F is overwritten with the UIName data. ---
$ws.Range("G2").Value2 = $ws.Range("F2").Value2
$ws.Range("G3").Value2 = $ws.Range("F3").Value2
$ws.Range("G4").Value2 = $ws.Range("F4").Value2

# Row 1 header for the (new) G column, copying F1's existing formatting
# (general alignment, same as the other header cells) so no new style is
# introduced.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value2 = "Description"

# --- Give column F the formatting used by the "type" column (C), i.e. the
#     general-alignment style, instead of the left-aligned one it had as the
#     old Description column. ---
$ws.Range("C2").Copy()
$ws.Range("F2").PasteSpecial(-4122)

$ws.Range("C3").Copy()
$ws.Range("F3").PasteSpecial(-4122)

$ws.Range("C4").Copy()
$ws.Range("F4").PasteSpecial(-4122)

# --- Now fill in the new UIName column values. ---
$ws.Range("F1").Value2 = "UIName"
$ws.Range("F2").Value2 = "武器UI名稱"
$ws.Range("F3").Value2 = "string"
$ws.Range("F4").Value2 = "Weapons Sprite Sheet_123"

# --- Split the old combined UIPath value into plain path + sprite name. ---
$ws.Range("E4").Value2 = "Assets/ArtResources/Weapons/Weapons Sprite Sheet.png"

# --- Size the new column G to roughly match the authored width (27.13 ch). ---
$ws.Columns("G:G").ColumnWidth = 26.3
